$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B through AB) must be swapped,
# while column A (the record index) stays untouched.
$pairs = @(
    @(26, 27),
    @(40, 41),
    @(60, 61),
    @(64, 65),
    @(70, 71),
    @(94, 95),
    @(124, 125),
    @(180, 181),
    @(182, 183),
    @(198, 199),
    @(208, 209),
    @(210, 211),
    @(229, 230)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
